$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet 1: Battery_Data
# Adds "upgrade 2" / "upgrade 3" rows for Nominal Capacity,
# Investment and Yearly O&M Cost, and zeroes out all the
# previously-computed battery figures (Salvage Value /
# Battery Replacement feature reset the cached results).
# ============================================================
$ws1 = $wb.Worksheets.Item("Battery_Data")

# First, propagate the bold/border label style (currently on A5:A8)
# down onto the new rows A9:A14 so every label cell matches.
$ws1.Range("A8").Copy($ws1.Range("A9:A14"))

$ws1.Range("A5").Value = "Nominal Capacity at upgrade 1"
$ws1.Range("B5").Value = 0
$ws1.Range("A6").Value = "Nominal Capacity at upgrade 2"
$ws1.Range("B6").Value = 0
$ws1.Range("A7").Value = "Nominal Capacity at upgrade 3"
$ws1.Range("B7").Value = 0

$ws1.Range("A8").Value = "Investment at upgrade 1"
$ws1.Range("B8").Value = 0
$ws1.Range("A9").Value = "Investment at upgrade 2"
$ws1.Range("B9").Value = 0
$ws1.Range("A10").Value = "Investment at upgrade 3"
$ws1.Range("B10").Value = 0

$ws1.Range("A11").Value = "Yearly O&M Cost at upgrade 1"
$ws1.Range("B11").Value = 0
$ws1.Range("A12").Value = "Yearly O&M Cost at upgrade 2"
$ws1.Range("B12").Value = 0
$ws1.Range("A13").Value = "Yearly O&M Cost at upgrade 3"
$ws1.Range("B13").Value = 0

$ws1.Range("A14").Value = "Total actualized Battery Reposition Cost"
$ws1.Range("B14").Value = 0

# ============================================================
# Sheet 2: Yearly BRC
# Only the Scenario 1 header plus "y = 1..3" rows remain; the
# rest (y = 4..20) are removed and the kept rows are reset to 0.
# ============================================================
$ws2 = $wb.Worksheets.Item("Yearly BRC")

$ws2.Range("B2").Value = 0
$ws2.Range("B3").Value = 0
$ws2.Range("B4").Value = 0

$ws2.Range("A5:B21").EntireRow.Delete()
